# Resize/reposition the block-diagram picture on Sheet1.
#
# The picture (the first/only shape on the sheet) is anchored with a
# two-cell anchor whose bottom-right corner moves from
#   col 12 / 229306 EMU, row 26 / 57800 EMU
# to
#   col 13 / 208557 EMU, row 29 / 28575 EMU
# i.e. the user dragged the picture's resize handle to make it noticeably
# bigger (while keeping its aspect ratio locked, per the picture's
# noChangeAspect lock) so it covers both the master and slave block
# diagrams being uploaded.
#
# The top-left corner of the picture does not move, so we only need to
# grow its Width/Height. Excel (and this COM surface) derives the
# anchor's bottom-right cell/offset from the shape's absolute
# Left/Top/Width/Height (in points), so we compute the Width/Height that
# reproduces the target bottom-right anchor exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$shp = $ws.Shapes.Item(1)

# New size, in points, chosen so the picture's right/bottom edge lands on
# the target anchor (col 13 / 208557 EMU, row 29 / 28575 EMU) while its
# top-left (Left/Top) is left untouched.
$shp.Width = 538.609311023622
$shp.Height = 409.5
